# "Change font single to resource"
#
# Typography sheet: fill in the "Wildcard Characters" column (G) for the
# two existing font rows (Medium/Small venus-rising-rg.ttf entries).
#
# Translation sheet: the text id in row 5 changes from the placeholder
# "SingleUseId3" to "ResourceId1", and the row 6 entry (SingleUseId4) is
# removed entirely (its cells are cleared).

$wb = $excel.ActiveWorkbook

$wsTypography = $wb.Worksheets.Item("Typography")
$wsTypography.Range("G4").Value = "> KM"
$wsTypography.Range("G5").Value = "%"

$wsTranslation = $wb.Worksheets.Item("Translation")
$wsTranslation.Range("B5").Value = "ResourceId1"
$wsTranslation.Range("B6:F6").ClearContents()
